$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.74203821656051
    "C2" = 0.691460055096419
    "D2" = 0.760299625468165
    "E2" = 0.700934579439252
    "F2" = 0.616724738675958

    "B3" = 0.455414012738854
    "C3" = 0.495867768595041
    "D3" = 0.49438202247191
    "E3" = 0.44392523364486
    "F3" = 0.338850174216028

    "B4" = 0.522292993630573
    "C4" = 0.534435261707989
    "D4" = 0.49063670411985
    "E4" = 0.47196261682243
    "F4" = 0.430313588850174

    "B5" = 0.331210191082803
    "C5" = 0.421487603305785
    "D5" = 0.322097378277154
    "E5" = 0.425233644859813
    "F5" = 0.376306620209059

    "B6" = 0.770700636942675
    "C6" = 0.774104683195592
    "D6" = 0.741573033707865
    "E6" = 0.831775700934579
    "F6" = 0.698606271777004
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
